$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.128.42"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.362.92"
$ws.Range("E3").Value = "  +0.54%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.26"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.41"
$ws.Range("E6").Value = "  -1.08%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  +5.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.105"
$ws.Range("E9").Value = "  +3.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.56"
$ws.Range("E10").Value = "  +2.55%  "

$ws.Range("E11").Value = "  -2.06%  "

$ws.Range("E12").Value = "  -1.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.06"
$ws.Range("E13").Value = "  +2.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.788.87"
$ws.Range("E14").Value = "  +0.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.073.29"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.378.95"
$ws.Range("E17").Value = "  +1.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.97"
$ws.Range("E18").Value = "  +3.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.32"
$ws.Range("E19").Value = "  +2.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "330.56"
$ws.Range("E20").Value = "  -1.18%  "

$ws.Range("E21").Value = "  +2.62%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "63.48"
$ws.Range("E23").Value = "  +2.90%  "

$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.25"
$ws.Range("E26").Value = "  -2.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.32"
$ws.Range("E27").Value = "  -6.30%  "

$ws.Range("E28").Value = "  -0.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.44"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  +1.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.39"
$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.995"
$ws.Range("E34").Value = "  -4.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.17"
$ws.Range("E36").Value = "  -0.14%  "

$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("E38").Value = "  -2.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.411"
$ws.Range("E39").Value = "  +8.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "142.73"
$ws.Range("E40").Value = "  -3.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  +2.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "288.37"
$ws.Range("E42").Value = "  +0.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0950"
$ws.Range("E43").Value = "  +2.46%  "

$ws.Range("E44").Value = "  +2.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.565"
$ws.Range("E45").Value = "  +0.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.86"
$ws.Range("E46").Value = "  -2.11%  "

$ws.Range("E47").Value = "  +2.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.392"
$ws.Range("E48").Value = "  +2.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.07"
$ws.Range("E49").Value = "  +0.01%  "

$ws.Range("E51").Value = "  +0.04%  "
